$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rodada 32 / 33 score updates (gols_mandante / gols_visitante) for games
# that were missing results before and now have been played.
$ws.Range("E26").Value = 1
$ws.Range("F26").Value = 1

$ws.Range("E31").Value = 1
$ws.Range("F31").Value = 0

$ws.Range("E92").Value = 2
$ws.Range("F92").Value = 2

$ws.Range("E93").Value = 0
$ws.Range("F93").Value = 0

$ws.Range("E94").Value = 1
$ws.Range("F94").Value = 2

$ws.Range("E95").Value = 0
$ws.Range("F95").Value = 3

$ws.Range("E96").Value = 2
$ws.Range("F96").Value = 0

$ws.Range("E97").Value = 1
$ws.Range("F97").Value = 0

$ws.Range("E98").Value = 0
$ws.Range("F98").Value = 2

$ws.Range("E99").Value = 3
$ws.Range("F99").Value = 0

$ws.Range("E100").Value = 0
$ws.Range("F100").Value = 1

$ws.Range("E101").Value = 1
$ws.Range("F101").Value = 0

$ws.Range("E102").Value = 2
$ws.Range("F102").Value = 0

$ws.Range("E103").Value = 1
$ws.Range("F103").Value = 0

$ws.Range("E104").Value = 1
$ws.Range("F104").Value = 2

$ws.Range("E105").Value = 0
$ws.Range("F105").Value = 0

$ws.Range("E106").Value = 2
$ws.Range("F106").Value = 1

$ws.Range("E107").Value = 3
$ws.Range("F107").Value = 0

$ws.Range("E108").Value = 0
$ws.Range("F108").Value = 0

$ws.Range("E109").Value = 2
$ws.Range("F109").Value = 1

$ws.Range("E110").Value = 2
$ws.Range("F110").Value = 1

$ws.Range("E111").Value = 0
$ws.Range("F111").Value = 0

$ws.Range("E112").Value = 2
$ws.Range("F112").Value = 2

$ws.Range("E113").Value = 0
$ws.Range("F113").Value = 0

$ws.Range("E114").Value = 1
$ws.Range("F114").Value = 2

# Row 115: game not played yet; highlight it in yellow instead of a score.
$ws.Range("E115:F115").Interior.Color = 65535

$ws.Range("E116").Value = 0
$ws.Range("F116").Value = 1

$ws.Range("E117").Value = 2
$ws.Range("F117").Value = 1

$ws.Range("E118").Value = 1
$ws.Range("F118").Value = 1

$ws.Range("E119").Value = 2
$ws.Range("F119").Value = 0

$ws.Range("E120").Value = 1
$ws.Range("F120").Value = 2

$ws.Range("E121").Value = 0
$ws.Range("F121").Value = 1

# Update the active selection to reflect where the user was working.
$ws.Range("F124").Select()
